$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the equipment/material values for the "Vertrieb" and "Reinigungskraft" rows
$ws.Range("C2").Value = "  PC-01 ,   AK-03"
$ws.Range("C4").Value = "RW-02"

# Update the active selection on the sheet
$ws.Range("F3").Select()
